$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Build the format (font size, border, wrap, vertical alignment) on Z2 first ---
$z2 = $ws.Range("Z2")
$z2.Font.Size = 7.5
$z2.Borders.LineStyle = 1
$z2.WrapText = $true
$z2.VerticalAlignment = -4108

# --- Copy that format onto AA2 and Z3 so they share the same style (no extra styles minted) ---
$z2.Copy()
$ws.Range("AA2").PasteSpecial(-4122)
$ws.Range("Z3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Set the new cell values (shared strings) ---
$ws.Range("Z2").Value = "Premium tolerance for vaious frequencies"
$ws.Range("AA2").Value = "Premium Tolerance limit for various Products"
$ws.Range("Z3").Value = "Limit of Tolerance in (whole number / decimal place)"

# --- Column widths for the two new columns ---
$ws.Columns.Item(26).ColumnWidth = 52
$ws.Columns.Item(27).ColumnWidth = 47.833333333333336

# --- Update view: scroll right and select the newly added range ---
$ws.Range("Z2:AA3").Select()
